# Auto-generated edit script: apply MIS data refresh (BAJAJ-PL base page update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Cells.Item(2, 1).Value2 = 'BKT0'
$ws.Cells.Item(2, 3).Value2 = 817120112.4299999
$ws.Cells.Item(2, 4).Value2 = 228
$ws.Cells.Item(2, 5).Value2 = 25
$ws.Cells.Item(2, 6).Value2 = 202
$ws.Cells.Item(2, 11).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 61842820.34
$ws.Cells.Item(2, 14).Value2 = 752412783.87
$ws.Cells.Item(2, 19).Value2 = 2864508.22
$ws.Cells.Item(2, 21).Value2 = 7.57
$ws.Cells.Item(2, 22).Value2 = 92.08
$ws.Cells.Item(2, 25).Value2 = 0.35
$ws.Cells.Item(2, 29).Value2 = 8915728
$ws.Cells.Item(2, 30).Value2 = 92.43000000000001
$ws.Cells.Item(2, 31).Value2 = 0.35

# Row 3
$ws.Cells.Item(3, 1).Value2 = 'BKT1'
$ws.Cells.Item(3, 3).Value2 = 288503013.31
$ws.Cells.Item(3, 4).Value2 = 199
$ws.Cells.Item(3, 5).Value2 = 27
$ws.Cells.Item(3, 6).Value2 = 125
$ws.Cells.Item(3, 8).Value2 = 44
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 13).Value2 = 39913237.29
$ws.Cells.Item(3, 14).Value2 = 194592171.52
$ws.Cells.Item(3, 16).Value2 = 50315521.08
$ws.Cells.Item(3, 19).Value2 = 3682083.42
$ws.Cells.Item(3, 21).Value2 = 13.83
$ws.Cells.Item(3, 22).Value2 = 67.45
$ws.Cells.Item(3, 25).Value2 = 1.28
$ws.Cells.Item(3, 26).Value2 = 17.44
$ws.Cells.Item(3, 29).Value2 = 5234507
$ws.Cells.Item(3, 30).Value2 = 86.17
$ws.Cells.Item(3, 31).Value2 = 18.72

# Row 4
$ws.Cells.Item(4, 1).Value2 = 'BKT1'
$ws.Cells.Item(4, 2).Value2 = 'RAJASTHAN'
$ws.Cells.Item(4, 3).Value2 = 9911622.050000001
$ws.Cells.Item(4, 4).Value2 = 15
$ws.Cells.Item(4, 5).Value2 = 15
$ws.Cells.Item(4, 11).Value2 = $null
$ws.Cells.Item(4, 13).Value2 = 9911622.050000001
$ws.Cells.Item(4, 19).Value2 = $null
$ws.Cells.Item(4, 21).Value2 = 100
$ws.Cells.Item(4, 25).Value2 = $null
$ws.Cells.Item(4, 29).Value2 = $null
$ws.Cells.Item(4, 30).Value2 = $null
$ws.Cells.Item(4, 31).Value2 = $null

# Row 5
$ws.Cells.Item(5, 1).Value2 = 'BKT2'
$ws.Cells.Item(5, 2).Value2 = 'RAJASTHAN'
$ws.Cells.Item(5, 3).Value2 = 1921240
$ws.Cells.Item(5, 5).Value2 = 2
$ws.Cells.Item(5, 11).Value2 = $null
$ws.Cells.Item(5, 13).Value2 = 1921240
$ws.Cells.Item(5, 19).Value2 = $null
$ws.Cells.Item(5, 21).Value2 = 100
$ws.Cells.Item(5, 25).Value2 = $null
$ws.Cells.Item(5, 29).Value2 = $null
$ws.Cells.Item(5, 30).Value2 = $null
$ws.Cells.Item(5, 31).Value2 = $null

# Row 6
$ws.Cells.Item(6, 1).Value2 = 'BKT3'
$ws.Cells.Item(6, 2).Value2 = 'RAJASTHAN'
$ws.Cells.Item(6, 3).Value2 = 1765526.64
$ws.Cells.Item(6, 4).Value2 = 4
$ws.Cells.Item(6, 5).Value2 = 4
$ws.Cells.Item(6, 13).Value2 = 1765526.64

# Row 7
$ws.Cells.Item(7, 1).Value2 = 'BKT4'
$ws.Cells.Item(7, 2).Value2 = 'RAJASTHAN'
$ws.Cells.Item(7, 3).Value2 = 489693
$ws.Cells.Item(7, 4).Value2 = 2
$ws.Cells.Item(7, 5).Value2 = 2
$ws.Cells.Item(7, 13).Value2 = 489693
$ws.Cells.Item(7, 21).Value2 = 100

# Row 8
$ws.Cells.Item(8, 1).Value2 = 'BKT5'
$ws.Cells.Item(8, 2).Value2 = 'RAJASTHAN'
$ws.Cells.Item(8, 3).Value2 = 2298733
$ws.Cells.Item(8, 4).Value2 = 2
$ws.Cells.Item(8, 5).Value2 = 2
$ws.Cells.Item(8, 13).Value2 = 2298733
$ws.Cells.Item(8, 21).Value2 = 100

# Row 9
$ws.Cells.Item(9, 1).Value2 = 'BKT7'
$ws.Cells.Item(9, 2).Value2 = 'RAJASTHAN'
$ws.Cells.Item(9, 3).Value2 = 1180938.14
$ws.Cells.Item(9, 4).Value2 = 1
$ws.Cells.Item(9, 5).Value2 = 1
$ws.Cells.Item(9, 13).Value2 = 1180938.14
$ws.Cells.Item(9, 21).Value2 = 100

